$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1780343333333333
$ws.Range("H2").Value = 0.534103
$ws.Range("I2").Value = 0.003649670474736916
$ws.Range("J2").Value = 0.003649670474736915
$ws.Range("Q2").Value = 0.01110845222833333
$ws.Range("R2").Value = 0.09997607005499999
$ws.Range("S2").Value = 0.003649670474736916
$ws.Range("T2").Value = 0.003649670474736915

# Row 3
$ws.Range("I3").Value = 0.09908483984804967
$ws.Range("J3").Value = 0.09908483984804965
$ws.Range("S3").Value = 0.09908483984804967
$ws.Range("T3").Value = 0.09908483984804965

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 35.63223
$ws.Range("H4").Value = 106.89669
$ws.Range("I4").Value = 0.7304540385283456
$ws.Range("J4").Value = 0.7304540385283456
$ws.Range("Q4").Value = 2.22327299085
$ws.Range("R4").Value = 20.00945691765
$ws.Range("S4").Value = 0.7304540385283456
$ws.Range("T4").Value = 0.7304540385283456

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5521946666666667
$ws.Range("H5").Value = 1.656584
$ws.Range("I5").Value = 0.011319887201011
$ws.Range("J5").Value = 0.011319887201011
$ws.Range("Q5").Value = 0.03445418622666667
$ws.Range("R5").Value = 0.31008767604
$ws.Range("S5").Value = 0.011319887201011
$ws.Range("T5").Value = 0.011319887201011

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 7.585023666666667
$ws.Range("H6").Value = 22.755071
$ws.Range("I6").Value = 0.155491563947857
$ws.Range("J6").Value = 0.1554915639478569
$ws.Range("Q6").Value = 0.4732675516816667
$ws.Range("R6").Value = 4.259407965135
$ws.Range("S6").Value = 0.155491563947857
$ws.Range("T6").Value = 0.1554915639478569
